# Add one entry into the mouse single cell dataset table (Mouse worksheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mouse")

# Rename the "Seq Depth" column header to "Seq Depth/cell" on the Mouse sheet.
$ws.Range("G1").Value = "Seq Depth/cell"

# Tidy up two existing values that drop the now-redundant "/cell" and "/ sample" suffixes.
$ws.Range("G7").Value = "0.43M UMI"
$ws.Range("G11").Value = "40,000 – 60,000 reads "

# Append the new Loo et al. row (row 12) to the table.
$ws.Range("A12").Value = "<a href=”https://www.nature.com/articles/s41467-018-08079-9#Bib1” target=”_blank”>Loo</a>"
$ws.Range("B12").Value = "Drop-seq"
$ws.Range("C12").Value = "3’"
$ws.Range("D12").Value = "E14.5, P0"
$ws.Range("E12").Value = "10,931 at E14.5; 7614 at P0"
$ws.Range("F12").Value = "Cortex"
$ws.Range("G12").Value = "12,000 reads"
$ws.Range("H12").Value = 1600

$ws.Range("G13").Select()
